$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": update match rows 2-4 with new fixtures/probabilities ---
$summary = $wb.Worksheets.Item("Summary")

# Row 2: Адмирал – Динамо Мн
$summary.Cells.Item(2, 1).Value = 1369
$summary.Cells.Item(2, 2).Value = 45991.41666666666
$summary.Cells.Item(2, 3).Value = 'Адмирал'
$summary.Cells.Item(2, 4).Value = 'Динамо Мн'
$summary.Cells.Item(2, 5).Value = 'Адмирал – Динамо Мн'
$summary.Cells.Item(2, 6).Value = 897831
$summary.Cells.Item(2, 7).Value = 'https://text.khl.ru/text/897831.html'
$summary.Cells.Item(2, 8).Value = 2.617412
$summary.Cells.Item(2, 9).Value = 4.5
$summary.Cells.Item(2, 10).Value = 7.117412
$summary.Cells.Item(2, 11).Value = 28.706351
$summary.Cells.Item(2, 12).Value = 36.917927
$summary.Cells.Item(2, 13).Value = 65.624278
$summary.Cells.Item(2, 14).Value = 0.284707
$summary.Cells.Item(2, 15).Value = 0.176781
$summary.Cells.Item(2, 16).Value = 0.53833
$summary.Cells.Item(2, 17).Value = 3.512382905934873
$summary.Cells.Item(2, 18).Value = 5.656716502339052
$summary.Cells.Item(2, 19).Value = 1.857596641465272
$summary.Cells.Item(2, 20).Value = 28.4707
$summary.Cells.Item(2, 21).Value = 17.6781
$summary.Cells.Item(2, 22).Value = 53.833
$summary.Cells.Item(2, 23).Value = 0.459734
$summary.Cells.Item(2, 24).Value = 0.540084
$summary.Cells.Item(2, 25).Value = 1.851563830811503
$summary.Cells.Item(2, 26).Value = 0.634992
$summary.Cells.Item(2, 27).Value = 0.364827
$summary.Cells.Item(2, 28).Value = 2.741025198244647
$summary.Cells.Item(2, 29).Value = 0.77787
$summary.Cells.Item(2, 30).Value = 0.221948
$summary.Cells.Item(2, 31).Value = 4.505559860868312
$summary.Cells.Item(2, 32).Value = 0.614955
$summary.Cells.Item(2, 33).Value = 0.385045
$summary.Cells.Item(2, 34).Value = 1.626135245668382
$summary.Cells.Item(2, 35).Value = 0.344695
$summary.Cells.Item(2, 36).Value = 0.655305
$summary.Cells.Item(2, 37).Value = 2.901115478901638
$summary.Cells.Item(2, 38).Value = 0.771042
$summary.Cells.Item(2, 39).Value = 0.228958
$summary.Cells.Item(2, 40).Value = 1.296946210452868
$summary.Cells.Item(2, 41).Value = 0.533521
$summary.Cells.Item(2, 42).Value = 0.466479
$summary.Cells.Item(2, 43).Value = 1.874340466448368
$summary.Cells.Item(2, 44).Value = 0.644371
$summary.Cells.Item(2, 45).Value = 1.551901001131336
$summary.Cells.Item(2, 46).Value = 0.850298
$summary.Cells.Item(2, 47).Value = 1.176058276039694

# Row 3: Амур – ХК Сочи
$summary.Cells.Item(3, 1).Value = 1369
$summary.Cells.Item(3, 2).Value = 45991.41666666666
$summary.Cells.Item(3, 3).Value = 'Амур'
$summary.Cells.Item(3, 4).Value = 'ХК Сочи'
$summary.Cells.Item(3, 5).Value = 'Амур – ХК Сочи'
$summary.Cells.Item(3, 6).Value = 897832
$summary.Cells.Item(3, 7).Value = 'https://text.khl.ru/text/897832.html'
$summary.Cells.Item(3, 8).Value = 1.225758
$summary.Cells.Item(3, 9).Value = 0.961538
$summary.Cells.Item(3, 10).Value = 2.187297
$summary.Cells.Item(3, 11).Value = 24.98031
$summary.Cells.Item(3, 12).Value = 25.615145
$summary.Cells.Item(3, 13).Value = 50.595455
$summary.Cells.Item(3, 14).Value = 0.827776
$summary.Cells.Item(3, 15).Value = 0.091762
$summary.Cells.Item(3, 16).Value = 0.07685
$summary.Cells.Item(3, 17).Value = 1.208056285758466
$summary.Cells.Item(3, 18).Value = 10.89775724155969
$summary.Cells.Item(3, 19).Value = 13.01236174365647
$summary.Cells.Item(3, 20).Value = 82.77759999999999
$summary.Cells.Item(3, 21).Value = 9.1762
$summary.Cells.Item(3, 22).Value = 7.685
$summary.Cells.Item(3, 23).Value = 0.350229
$summary.Cells.Item(3, 24).Value = 0.646159
$summary.Cells.Item(3, 25).Value = 1.547606703613197
$summary.Cells.Item(3, 26).Value = 0.520867
$summary.Cells.Item(3, 27).Value = 0.475522
$summary.Cells.Item(3, 28).Value = 2.10295212419194
$summary.Cells.Item(3, 29).Value = 0.678623
$summary.Cells.Item(3, 30).Value = 0.317765
$summary.Cells.Item(3, 31).Value = 3.146979686246125
$summary.Cells.Item(3, 32).Value = 0.917635
$summary.Cells.Item(3, 33).Value = 0.08236499999999999
$summary.Cells.Item(3, 34).Value = 1.089757910280231
$summary.Cells.Item(3, 35).Value = 0.780622
$summary.Cells.Item(3, 36).Value = 0.219378
$summary.Cells.Item(3, 37).Value = 1.281029742948572
$summary.Cells.Item(3, 38).Value = 0.413325
$summary.Cells.Item(3, 39).Value = 0.5866749999999999
$summary.Cells.Item(3, 40).Value = 2.419403617008407
$summary.Cells.Item(3, 41).Value = 0.170132
$summary.Cells.Item(3, 42).Value = 0.8298680000000001
$summary.Cells.Item(3, 43).Value = 5.877789010885665
$summary.Cells.Item(3, 44).Value = 0.967298
$summary.Cells.Item(3, 45).Value = 1.033807575328389
$summary.Cells.Item(3, 46).Value = 0.308085
$summary.Cells.Item(3, 47).Value = 3.245857474398299

# Row 4: Ак Барс – Драконы
$summary.Cells.Item(4, 1).Value = 1369
$summary.Cells.Item(4, 2).Value = 45991.70833333334
$summary.Cells.Item(4, 3).Value = 'Ак Барс'
$summary.Cells.Item(4, 4).Value = 'Драконы'
$summary.Cells.Item(4, 5).Value = 'Ак Барс – Драконы'
$summary.Cells.Item(4, 6).Value = 897833
$summary.Cells.Item(4, 7).Value = 'https://text.khl.ru/text/897833.html'
$summary.Cells.Item(4, 8).Value = 3.055625
$summary.Cells.Item(4, 9).Value = 3.676454
$summary.Cells.Item(4, 10).Value = 6.732079
$summary.Cells.Item(4, 11).Value = 33.778683
$summary.Cells.Item(4, 12).Value = 31.087527
$summary.Cells.Item(4, 13).Value = 64.86621100000001
$summary.Cells.Item(4, 14).Value = 0.49328
$summary.Cells.Item(4, 15).Value = 0.140849
$summary.Cells.Item(4, 16).Value = 0.359373
$summary.Cells.Item(4, 17).Value = 2.027246188777165
$summary.Cells.Item(4, 18).Value = 7.099801915526557
$summary.Cells.Item(4, 19).Value = 2.782624181560663
$summary.Cells.Item(4, 20).Value = 49.328
$summary.Cells.Item(4, 21).Value = 14.0849
$summary.Cells.Item(4, 22).Value = 35.9373
$summary.Cells.Item(4, 23).Value = 0.096294
$summary.Cells.Item(4, 24).Value = 0.897208
$summary.Cells.Item(4, 25).Value = 1.11456875105884
$summary.Cells.Item(4, 26).Value = 0.185875
$summary.Cells.Item(4, 27).Value = 0.807627
$summary.Cells.Item(4, 28).Value = 1.238195355033945
$summary.Cells.Item(4, 29).Value = 0.3062
$summary.Cells.Item(4, 30).Value = 0.687302
$summary.Cells.Item(4, 31).Value = 1.454964484316938
$summary.Cells.Item(4, 32).Value = 0.926781
$summary.Cells.Item(4, 33).Value = 0.07321900000000001
$summary.Cells.Item(4, 34).Value = 1.079003561790757
$summary.Cells.Item(4, 35).Value = 0.7998459999999999
$summary.Cells.Item(4, 36).Value = 0.200154
$summary.Cells.Item(4, 37).Value = 1.250240671329231
$summary.Cells.Item(4, 38).Value = 0.891025
$summary.Cells.Item(4, 39).Value = 0.108975
$summary.Cells.Item(4, 40).Value = 1.122302965685587
$summary.Cells.Item(4, 41).Value = 0.7280799999999999
$summary.Cells.Item(4, 42).Value = 0.27192
$summary.Cells.Item(4, 43).Value = 1.373475442259092
$summary.Cells.Item(4, 44).Value = 0.758018
$summary.Cells.Item(4, 45).Value = 1.319229886361538
$summary.Cells.Item(4, 46).Value = 0.6403759999999999
$summary.Cells.Item(4, 47).Value = 1.561582570239984

# --- Sheet "Cards_telegram": update match rows 2-4 with new fixtures/card text ---
$cards = $wb.Worksheets.Item("Cards_telegram")

# Row 2: Адмирал – Динамо Мн
$cards.Cells.Item(2, 1).Value = 45991.41666666666
$cards.Cells.Item(2, 2).Value = 'Адмирал – Динамо Мн'
$cards.Cells.Item(2, 3).Value = @'
КХЛ • Регулярный чемпионат • 30.11.2025

Адмирал – Динамо Мн

Ожидания модели (60’):
• Голы: λ_total ≈ 4.89 (2.08 : 2.81)
• Броски: SOG λ ≈ 66 (29 : 37)

Исход (60’), честные кф:
• П1: 28.5%  (Kмод 3.51)
• Х:  17.7%  (Kмод 5.66)
• П2: 53.8%  (Kмод 1.86)

Тоталы голов:
• ТМ 4.5: 46.0%  (Kмод 2.18)
• ТБ 4.5: 54.0%  (Kмод 1.85)

• ТМ 5.5: 63.5%  (Kмод 1.57)
• ТБ 5.5: 36.5%  (Kмод 2.74)

• ТМ 6.5: 77.8%  (Kмод 1.29)
• ТБ 6.5: 22.2%  (Kмод 4.51)

Индивидуальные тоталы:
• Адмирал ИТБ 1.5: 61.5% (Kмод 1.63)
• Адмирал ИТБ 2.5: 34.5% (Kмод 2.90)
• Динамо Мн ИТБ 1.5: 77.1% (Kмод 1.30)
• Динамо Мн ИТБ 2.5: 53.4% (Kмод 1.87)

Фора +1.5:
• Адмирал +1.5: 64.4% (Kмод 1.55)
• Динамо Мн +1.5: 85.0% (Kмод 1.18)
'@

# Row 3: Амур – ХК Сочи
$cards.Cells.Item(3, 1).Value = 45991.41666666666
$cards.Cells.Item(3, 2).Value = 'Амур – ХК Сочи'
$cards.Cells.Item(3, 3).Value = @'
КХЛ • Регулярный чемпионат • 30.11.2025

Амур – ХК Сочи

Ожидания модели (60’):
• Голы: λ_total ≈ 5.55 (4.13 : 1.41)
• Броски: SOG λ ≈ 51 (25 : 26)

Исход (60’), честные кф:
• П1: 82.8%  (Kмод 1.21)
• Х:  9.2%  (Kмод 10.90)
• П2: 7.7%  (Kмод 13.01)

Тоталы голов:
• ТМ 4.5: 35.0%  (Kмод 2.86)
• ТБ 4.5: 64.6%  (Kмод 1.55)

• ТМ 5.5: 52.1%  (Kмод 1.92)
• ТБ 5.5: 47.6%  (Kмод 2.10)

• ТМ 6.5: 67.9%  (Kмод 1.47)
• ТБ 6.5: 31.8%  (Kмод 3.15)

Индивидуальные тоталы:
• Амур ИТБ 1.5: 91.8% (Kмод 1.09)
• Амур ИТБ 2.5: 78.1% (Kмод 1.28)
• ХК Сочи ИТБ 1.5: 41.3% (Kмод 2.42)
• ХК Сочи ИТБ 2.5: 17.0% (Kмод 5.88)

Фора +1.5:
• Амур +1.5: 96.7% (Kмод 1.03)
• ХК Сочи +1.5: 30.8% (Kмод 3.25)
'@

# Row 4: Ак Барс – Драконы
$cards.Cells.Item(4, 1).Value = 45991.70833333334
$cards.Cells.Item(4, 2).Value = 'Ак Барс – Драконы'
$cards.Cells.Item(4, 3).Value = @'
КХЛ • Регулярный чемпионат • 30.11.2025

Ак Барс – Драконы

Ожидания модели (60’):
• Голы: λ_total ≈ 8.06 (4.28 : 3.78)
• Броски: SOG λ ≈ 65 (34 : 31)

Исход (60’), честные кф:
• П1: 49.3%  (Kмод 2.03)
• Х:  14.1%  (Kмод 7.10)
• П2: 35.9%  (Kмод 2.78)

Тоталы голов:
• ТМ 4.5: 9.6%  (Kмод 10.38)
• ТБ 4.5: 89.7%  (Kмод 1.11)

• ТМ 5.5: 18.6%  (Kмод 5.38)
• ТБ 5.5: 80.8%  (Kмод 1.24)

• ТМ 6.5: 30.6%  (Kмод 3.27)
• ТБ 6.5: 68.7%  (Kмод 1.45)

Индивидуальные тоталы:
• Ак Барс ИТБ 1.5: 92.7% (Kмод 1.08)
• Ак Барс ИТБ 2.5: 80.0% (Kмод 1.25)
• Драконы ИТБ 1.5: 89.1% (Kмод 1.12)
• Драконы ИТБ 2.5: 72.8% (Kмод 1.37)

Фора +1.5:
• Ак Барс +1.5: 75.8% (Kмод 1.32)
• Драконы +1.5: 64.0% (Kмод 1.56)
'@

